# "Update the draft and fix almost all Missed rate calculation result."
#
# 1) Drop the orphaned external reference (xl/externalLinks/externalLink1.xml,
#    and the <externalReferences> node in workbook.xml) - nothing in the
#    workbook actually formula-references it, it's just a stale link to
#    er_memory.xlsx left over from a copy/paste.
# 2) Recompute the "GR miss risk" / "PR optimal miss risk" rows (16 & 17) on
#    the "walking" sheet - these are plain cached numbers (no formulas in
#    this workbook), so we overwrite them with the corrected values.
# 3) Carry forward the cursor/selection that was left in each sheet when the
#    author saved the file.

$wb = $excel.ActiveWorkbook

# --- 1) Remove the stale external link to er_memory.xlsx -------------------
foreach ($link in @($wb.LinkSources())) {
    $wb.BreakLink($link, 1)
}

# --- 2) Fix the Missed-rate ("miss risk") calculation results --------------
$walking = $wb.Worksheets.Item("walking")

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

# Row 16 = "GR miss risk"
$grMissRisk = @(
    61.895099238493998,
    64.855108869195405,
    64.56569456820435,
    63.960593231504333,
    63.155326380440492,
    62.220579348256223,
    61.587561475131672,
    60.522799305501621,
    59.530757606572934,
    58.65542647664109
)

# Row 17 = "PR optimal miss risk"
$prOptimalMissRisk = @(
    2.06353112128674,
    2.4918225452438403,
    2.8153417648856101,
    3.0881471354135801,
    3.3579077141734199,
    3.6367871630851201,
    3.9427732146170098,
    4.2576385151629799,
    4.6237726900929204,
    5.01129407199425
)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $walking.Range($cols[$i] + "16").Value = $grMissRisk[$i]
    $walking.Range($cols[$i] + "17").Value = $prOptimalMissRisk[$i]
}

# --- 3) Restore each sheet's last-used selection ----------------------------
$walking.Range("E10").Select()

$weekday = $wb.Worksheets.Item("weekday")
$weekday.Range("B9").Select()

$erMemory = $wb.Worksheets.Item("ER_memory_and_function")
$erMemory.Range("F1").Select()
